$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 37
$ws.Range("E8").Value = 39
$ws.Range("E9").Value = 58
$ws.Range("E24").Value = 26
$ws.Range("E27").Value = 8
$ws.Range("E31").Value = 55
$ws.Range("E36").Value = 42
$ws.Range("E39").Value = 63
$ws.Range("E43").Value = 18
$ws.Range("E46").Value = 88
$ws.Range("E47").Value = 38
